# edit.ps1 - PowerShell-style PowerPoint COM-interop script
#
# Two changes, matching the supplied OOXML diff:
#
#  1) The table on slide 16 (the "Cash flow" plenary recap table) switches
#     from the deck's custom table style ("Table_0") to the built-in table
#     style {80C48809-8081-4FAB-95E0-FF283BAF0775}.
#
#  2) The presentation's theme colour scheme is repainted from the
#     "Integral" palette to the stock "Office Theme" palette (dk1/lt1/dk2/
#     lt2/accent1-6/hlink/folHlink). fontScheme and fmtScheme are identical
#     between the two themes, so only the twelve colour slots need updating.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------

$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{80C48809-8081-4FAB-95E0-FF283BAF0775}")
    }
}

# --- 2) Theme colours: Integral -> Office Theme ---------------------------

function Set-ThemeColor {
    param($colorScheme, [int]$index, [string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

# index : scheme slot   : Integral (before) -> Office Theme (after)
Set-ThemeColor $themeColors 1  "000000"   # dk1
Set-ThemeColor $themeColors 2  "FFFFFF"   # lt1
Set-ThemeColor $themeColors 3  "44546A"   # dk2
Set-ThemeColor $themeColors 4  "E7E6E6"   # lt2
Set-ThemeColor $themeColors 5  "5B9BD5"   # accent1
Set-ThemeColor $themeColors 6  "ED7D31"   # accent2
Set-ThemeColor $themeColors 7  "A5A5A5"   # accent3
Set-ThemeColor $themeColors 8  "FFC000"   # accent4
Set-ThemeColor $themeColors 9  "4472C4"   # accent5
Set-ThemeColor $themeColors 10 "70AD47"   # accent6
Set-ThemeColor $themeColors 11 "0563C1"   # hlink
Set-ThemeColor $themeColors 12 "954F72"   # folHlink
